$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 51: PingPongDelay / Left Delay Time ---
$ws.Range("A51").Value = "PingPongDelay"
$ws.Range("D51").Value = "L-R Delay Time"
$ws.Range("B51").Value = "Left Delay Time"
$ws.Range("C51").Value = "delayTimeL"
$ws.Range("K51").Value = "delaySecL"
$ws.Range("F51").Value = "sec"
$ws.Range("G51").Value = 0.01
$ws.Range("H51").Value = 2
$ws.Range("I51").Value = 0.5
$ws.Range("J51").Value = 0.01

# --- Row 52: PingPongDelay / Right Delay Time ---
$ws.Range("A52").Value = "PingPongDelay"
$ws.Range("B52").Value = "Right Delay Time"
$ws.Range("C52").Value = "delayTimeR"
$ws.Range("D52").Value = "R-L Delay Time"
$ws.Range("K52").Value = "delaySecR"
$ws.Range("F52").Value = "sec"
$ws.Range("G52").Value = 0.01
$ws.Range("H52").Value = 2
$ws.Range("I52").Value = 0.5
$ws.Range("J52").Value = 0.01

# --- Row 53: PingPongDelay / Feedback ---
$ws.Range("A53").Value = "PingPongDelay"
$ws.Range("B53").Value = "Feedback"
$ws.Range("C53").Value = "feedback"
$ws.Range("D53").Value = "Feedback"
$ws.Range("K53").Value = "feedback"
$ws.Range("F53").Value = "%"
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 99.5
$ws.Range("I53").Value = 75
$ws.Range("J53").Value = 0.5

# --- Row 54: PingPongDelay / Wet Mix Level ---
$ws.Range("A54").Value = "PingPongDelay"
$ws.Range("B54").Value = "Wet Mix Level"
$ws.Range("C54").Value = "wetLevel"
$ws.Range("D54").Value = "Delay Mix Level"
$ws.Range("K54").Value = "wetLevel"
$ws.Range("F54").Value = "%"
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 100
$ws.Range("I54").Value = 50
$ws.Range("J54").Value = 1

# Update selection to match the new active cell after editing (B57)
$ws.Range("B57").Select()
